$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 0.37990752884762458
$ws.Cells.Item(2, 1).Value = -0.0099999993713915103
$ws.Cells.Item(3, 1).Value = -0.0089999993580871518
$ws.Cells.Item(4, 1).Value = -0.011999999829235719
$ws.Cells.Item(5, 1).Value = -0.0059999993625829973
$ws.Cells.Item(6, 1).Value = -0.0059999993415367214
$ws.Cells.Item(7, 1).Value = -0.019999999228947019
$ws.Cells.Item(8, 1).Value = -0.019999999225461806
$ws.Cells.Item(9, 1).Value = -0.0059999993344366231
$ws.Cells.Item(10, 1).Value = -0.0059999993333121893
$ws.Cells.Item(11, 1).Value = -0.0044999993451853015
$ws.Cells.Item(12, 1).Value = -0.0059999993333739177
$ws.Cells.Item(13, 1).Value = -0.0059999993355477343
$ws.Cells.Item(14, 1).Value = -0.011999999289254326
$ws.Cells.Item(15, 1).Value = 0.048903961416226771
$ws.Cells.Item(16, 1).Value = -0.0059999993368755611
$ws.Cells.Item(17, 1).Value = -0.0059999993337136459
$ws.Cells.Item(18, 1).Value = -0.009861174282823626
$ws.Cells.Item(19, 1).Value = -0.0089999993775564668
$ws.Cells.Item(20, 1).Value = -0.060609865892079284
$ws.Cells.Item(21, 1).Value = -0.0089999993582621229
$ws.Cells.Item(22, 1).Value = -0.0089999993574343407
$ws.Cells.Item(23, 1).Value = -0.0089999993524871869
$ws.Cells.Item(24, 1).Value = -0.041999999079017414
$ws.Cells.Item(25, 1).Value = -0.041999999073254024
$ws.Cells.Item(26, 1).Value = -0.005999999339135087
$ws.Cells.Item(27, 1).Value = -0.0049696746610341513
$ws.Cells.Item(28, 1).Value = -0.0059999993207346947
$ws.Cells.Item(29, 1).Value = -0.011999999262741312
$ws.Cells.Item(30, 1).Value = -0.019999999193859086
$ws.Cells.Item(31, 1).Value = -0.014999999226020577
$ws.Cells.Item(32, 1).Value = -0.023234921353948224
$ws.Cells.Item(33, 1).Value = -0.0059999992958399417
